# Mise à jour de l'application
# Add the "N3J2" matchday results (columns DC:DF) for players who took part.
# Layout per match block: <min> | T/R | But (goal) | Passe D (assist)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("DC2").Value = 90
$ws.Range("DD2").Value = "T"

$ws.Range("DC5").Value = 90
$ws.Range("DD5").Value = "T"

$ws.Range("DC9").Value = 90
$ws.Range("DD9").Value = "T"

$ws.Range("DC11").Value = 78
$ws.Range("DD11").Value = "T"

$ws.Range("DC12").Value = 49
$ws.Range("DD12").Value = "T"

$ws.Range("DC13").Value = 65
$ws.Range("DD13").Value = "T"

$ws.Range("DC14").Value = 25
$ws.Range("DD14").Value = "R"

$ws.Range("DC15").Value = 72
$ws.Range("DD15").Value = "T"
$ws.Range("DE15").Value = 1

$ws.Range("DC16").Value = 90
$ws.Range("DD16").Value = "T"
$ws.Range("DF16").Value = 1

$ws.Range("DC18").Value = 18
$ws.Range("DD18").Value = "R"

$ws.Range("DC20").Value = 78
$ws.Range("DD20").Value = "T"

$ws.Range("DC22").Value = 90
$ws.Range("DD22").Value = "T"

$ws.Range("DC24").Value = 90
$ws.Range("DD24").Value = "T"

$ws.Range("DC25").Value = 12
$ws.Range("DD25").Value = "R"

$ws.Range("DC26").Value = 12
$ws.Range("DD26").Value = "R"

# Reflect the view state change recorded in the saved workbook: the
# active selection moved to DH16 (the frozen pane stays split after
# column A; only the active cell is something we can reliably replay).
$ws.Range("DH16").Select()
